$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 (2020-03-06): new positive cases corrected 95 -> 96
$ws.Range("C24").Value = 96

# Row 389 (2021-03-20): new positive cases corrected 60 -> 59
$ws.Range("C389").Value = 59

# Row 401 (2021-04-01): new positive cases corrected 81 -> 82
$ws.Range("C401").Value = 82

# Row 403 (2021-04-03): new positive cases corrected 42 -> 41; one new hospital death recorded
$ws.Range("C403").Value = 41
$ws.Range("L403").Value = 1

# Row 405 (2021-04-05): new positive cases corrected 43 -> 57
$ws.Range("C405").Value = 57

# Row 406 (2021-04-06): new positive cases corrected 9 -> 85
$ws.Range("C406").Value = 85

# Row 407 (2021-04-07): fill in the day's data (previously blank placeholder row)
$ws.Range("C407").Value = 13
$ws.Range("E407").Value = 8
$ws.Range("F407").Value = 6
$ws.Range("G407").Value = 31
$ws.Range("L407").Value = 0
$ws.Range("M407").Value = 0
